# Spring Project Report - add "17.04.2023" status update section (rows 10-18)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B is widened to fit the long URLs now stored in it ---
$ws.Columns.Item(2).ColumnWidth = 118.140625
$ws.Columns.Item(3).ColumnWidth = 10.85546875
$ws.Columns.Item(4).ColumnWidth = 10.85546875

# --- existing header row 2 ("Reading" / 13.04.2023 block) gets centred+middle
#     vertical alignment now that the row is shorter, and its wrap height drops
$ws.Range("A2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A2").VerticalAlignment = -4108     # xlCenter
$ws.Rows.Item(2).RowHeight = 30.75

# --- existing "Reading" detail rows (3-7) get vertical centring added ---
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A4").HorizontalAlignment = -4108
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A5").HorizontalAlignment = -4108
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A6").HorizontalAlignment = -4108
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("A7").VerticalAlignment = -4108

# ======================================================================
# New status block: 17.04.2023
# ======================================================================

$ws.Range("B10").Value = "17.04.2023"

$ws.Range("A11").Value = "Reading"
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("A11").VerticalAlignment = -4108
$ws.Range("B11").Value = "Chapter 3. Forces and Equations of Motion, Section 3.7 Hook's Law and Simple"
$ws.Range("B11").Font.Color = 4473924

$ws.Range("A12").Value = "Youtoube videos"
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("A12").VerticalAlignment = -4108
$ws.Range("B12").Value = "https://www.youtube.com/watch?v=g550H4e5FCY&pp=ygUeRm9yY2VzIGFuZCBFcXVhdGlvbnMgb2YgTW90aW9u"

$ws.Range("A13").HorizontalAlignment = -4108
$ws.Range("A13").VerticalAlignment = -4108
$ws.Range("B13").Value = "https://www.youtube.com/watch?v=UNuRhIHthhY"

$ws.Range("A14").HorizontalAlignment = -4108
$ws.Range("A14").VerticalAlignment = -4108
$ws.Range("B14").Value = "https://www.youtube.com/watch?v=-_l_YDA6au8&pp=ygUqNyBIb29r4oCZcyBMYXcgYW5kIFNpbXBsZSBIYXJtb25pYyBNb3Rpb24s"

$ws.Range("A15").HorizontalAlignment = -4108
$ws.Range("A15").VerticalAlignment = -4108
$ws.Range("B15").Value = "https://www.youtube.com/watch?v=_Gnke2x3vT8&pp=ygUqNyBIb29r4oCZcyBMYXcgYW5kIFNpbXBsZSBIYXJtb25pYyBNb3Rpb24s"

$ws.Range("A16").Value = "Research"
$ws.Range("B16").Value = "hooke-s-law-and-simple-harmonic-motion"

$ws.Range("B17").Value = "A_Harmonic_Oscillator_Obeys_Hooke's_Law"

$ws.Range("A18").Value = "Task 2"
$ws.Range("B18").Value = "Completed :)"

# --- merges for the grouped label cells ---
$ws.Range("A12:A15").Merge()
$ws.Range("A16:A17").Merge()

# --- hyperlinks (order matches the relationship ids in the target file) ---
$ws.Hyperlinks.Add($ws.Range("B13"), "https://www.youtube.com/watch?v=UNuRhIHthhY")
$ws.Hyperlinks.Add($ws.Range("B14"), "https://www.youtube.com/watch?v=-_l_YDA6au8&pp=ygUqNyBIb29r4oCZcyBMYXcgYW5kIFNpbXBsZSBIYXJtb25pYyBNb3Rpb24s")
$ws.Hyperlinks.Add($ws.Range("B15"), "https://www.youtube.com/watch?v=_Gnke2x3vT8&pp=ygUqNyBIb29r4oCZcyBMYXcgYW5kIFNpbXBsZSBIYXJtb25pYyBNb3Rpb24s")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://library.fiveable.me/key-terms/ap-physics-1/hooke-s-law-and-simple-harmonic-motion")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://phys.libretexts.org/Bookshelves/University_Physics/University_Physics_(OpenStax)/Book%3A_University_Physics_I_-_Mechanics_Sound_Oscillations_and_Waves_(OpenStax)/15%3A_Oscillations/15.02%3A_Simple_Harmonic_Motion")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.youtube.com/watch?v=g550H4e5FCY&pp=ygUeRm9yY2VzIGFuZCBFcXVhdGlvbnMgb2YgTW90aW9u")

# Re-apply the plain "Hyperlink" look-alike style (same visual style already
# used by B3:B7) to the cells the Hyperlinks.Add() calls above just touched,
# so they line up with the rest of the sheet instead of minting a new style.
$ws.Range("B12").Style = "Hyperlink"
$ws.Range("B13").Style = "Hyperlink"
$ws.Range("B14").Style = "Hyperlink"
$ws.Range("B15").Style = "Hyperlink"
$ws.Range("B16").Style = "Hyperlink"
$ws.Range("B17").Style = "Hyperlink"

Write-Host "edit applied"
